# Weekly update: a new week of "Apio" (Vega Central Mapocho de Santiago)
# price data is inserted at the top of the data block (rows 149:150),
# pushing the existing history down by two rows. The two oldest rows
# that fall off the bottom of the original range re-appear as new rows
# 155:156 (dimension grows from A1:R154 to A1:R156).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new week: push rows 151:154 down to 153:156.
$ws.Rows("151:152").Insert()

# The new week re-uses last week's figures (rows 149:150) as its
# starting values, which then get the correct newer date below; the
# vacated rows 151:152 are filled from the (now shifted) previous week.
$ws.Range("A149:R150").Copy($ws.Range("A151:R152"))

# Stamp the new week's date onto the top two rows.
$ws.Range("D149").Value = 44516
$ws.Range("D150").Value = 44516
